$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1169995834814548, 0.002658071450198252, 3.223369029078222, 13.86384647080068, 17.20687315481056)
    3  = @(0.00009552326474482342, 0.002658071450198252, 0.7210945179870265, 13.86384647080068, 14.58769458350265)
    4  = @(0.2881169905109251, 0.002658071450198252, 0.1496068669990043, 0.5333859586016987, 0.9737678875618263)
    5  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    6  = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.038307959104277)
    7  = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 2.213936997104367)
    8  = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    9  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    10 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    11 = @(0.2881169905109251, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.276001064745008)
    12 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    13 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    14 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    15 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 24.14949828602258)
    16 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 19.48425592650926)
    17 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.178645819794754)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
